$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H55").Value = 347.4737
$ws.Range("I55").Value = 318.125
$ws.Range("J55").Value = 368.81818
$ws.Range("K55").Value = 318.125
$ws.Range("L55").Value = 368.81818
$ws.Range("M55").Value = -104.125
$ws.Range("N55").Value = -796.81818

$ws.Range("H62").Value = 3710.8125
$ws.Range("I62").Value = 3239.1667
$ws.Range("J62").Value = 5125.75
$ws.Range("K62").Value = 3239.1667
$ws.Range("L62").Value = 5125.75
$ws.Range("M62").Value = -2615.1667
$ws.Range("N62").Value = -6373.75

$ws.Range("H65").Value = 3710.8125
$ws.Range("I65").Value = 3239.1667
$ws.Range("J65").Value = 5125.75
$ws.Range("K65").Value = 16195.8335
$ws.Range("L65").Value = 25628.75
$ws.Range("M65").Value = -13075.8335
$ws.Range("N65").Value = -31868.75

$ws.Range("H92").Value = 612.53845
$ws.Range("I92").Value = 635.08
$ws.Range("K92").Value = 635.08
$ws.Range("M92").Value = 612.92

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 254.25
$ws.Range("I4").Value = 249.33333
$ws.Range("K4").Value = 249.33333
$ws.Range("M4").Value = -133.33333

$ws.Range("H32").Value = 2087.0652
$ws.Range("I32").Value = 2123.0698
$ws.Range("K32").Value = 2123.0698
$ws.Range("M32").Value = -1836.0698

$ws.Range("H45").Value = 15175809
$ws.Range("I45").Value = 34939.047
$ws.Range("K45").Value = 34939.047
$ws.Range("M45").Value = -34562.047

$ws.Range("H122").Value = 18869448
$ws.Range("J122").Value = 2119.6667
$ws.Range("L122").Value = 6359.000100000001
$ws.Range("N122").Value = -11259.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8259.75
$ws.Range("I132").Value = 9896.691999999999
$ws.Range("J132").Value = 1166.3334
$ws.Range("K132").Value = 29690.076
$ws.Range("L132").Value = 3499.0002
$ws.Range("M132").Value = -27160.076
$ws.Range("N132").Value = -8559.0002

$ws.Range("H138").Value = 4975.7715
$ws.Range("J138").Value = 5053.3228
$ws.Range("L138").Value = 15159.9684
$ws.Range("N138").Value = -25439.9684

$ws.Range("H139").Value = 69981.3
$ws.Range("J139").Value = 69998.71000000001
$ws.Range("L139").Value = 69998.71000000001
$ws.Range("N139").Value = -80278.71000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 523.8
$ws.Range("I80").Value = 284
$ws.Range("J80").Value = 583.75
$ws.Range("K80").Value = 284
$ws.Range("L80").Value = 583.75
$ws.Range("M80").Value = 714
$ws.Range("N80").Value = -2579.75

$ws.Range("H83").Value = 523.8
$ws.Range("I83").Value = 284
$ws.Range("J83").Value = 583.75
$ws.Range("K83").Value = 1420
$ws.Range("L83").Value = 2918.75
$ws.Range("M83").Value = 3572
$ws.Range("N83").Value = -12902.75

$ws.Range("H97").Value = 23426.637
$ws.Range("I97").Value = 2798
$ws.Range("J97").Value = 25489.5
$ws.Range("K97").Value = 2798
$ws.Range("L97").Value = 25489.5
$ws.Range("M97").Value = -1807
$ws.Range("N97").Value = -27471.5

$ws.Range("H99").Value = 3779.5334
$ws.Range("I99").Value = 3526.6365
$ws.Range("K99").Value = 3526.6365
$ws.Range("M99").Value = -2028.6365

$ws.Range("H105").Value = 18573240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3210.5762
$ws.Range("I31").Value = 1622.5385
$ws.Range("J31").Value = 3659.3696
$ws.Range("K31").Value = 1622.5385
$ws.Range("L31").Value = 3659.3696
$ws.Range("M31").Value = -1327.5385
$ws.Range("N31").Value = -4249.3696

$ws.Range("H34").Value = 3210.5762
$ws.Range("I34").Value = 1622.5385
$ws.Range("J34").Value = 3659.3696
$ws.Range("K34").Value = 1622.5385
$ws.Range("L34").Value = 3659.3696
$ws.Range("M34").Value = -1420.5385
$ws.Range("N34").Value = -4063.3696

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2389
$ws.Range("I134").Value = 2389
$ws.Range("K134").Value = 7167
$ws.Range("M134").Value = -4632

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2307.8948
$ws.Range("J107").Value = 2810.8667
$ws.Range("L107").Value = 8432.6001
$ws.Range("N107").Value = -12272.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3208.4473
$ws.Range("I132").Value = 2752.2083
$ws.Range("J132").Value = 3990.5715
$ws.Range("K132").Value = 8256.624899999999
$ws.Range("L132").Value = 11971.7145
$ws.Range("M132").Value = -5726.624899999999
$ws.Range("N132").Value = -17031.7145

$ws.Range("H134").Value = 4001.4285
$ws.Range("I134").Value = 4173.2173
$ws.Range("K134").Value = 12519.6519
$ws.Range("M134").Value = -9984.651900000001

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 83336230
$ws.Range("I113").Value = 90911930
$ws.Range("K113").Value = 90911930
$ws.Range("M113").Value = -90909760

$ws.Range("H114").Value = 12000
$ws.Range("J114").Value = 12000
$ws.Range("L114").Value = 12000
$ws.Range("N114").Value = -20678

$ws.Range("H126").Value = 41670564
$ws.Range("I126").Value = 45457932
$ws.Range("K126").Value = 136373796
$ws.Range("M126").Value = -136371326

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4658.6216
$ws.Range("J132").Value = 4627.613
$ws.Range("L132").Value = 41648.517
$ws.Range("N132").Value = -46708.517

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2974.5
$ws.Range("I7").Value = 2999.3333
$ws.Range("J7").Value = 2900
$ws.Range("K7").Value = 2999.3333
$ws.Range("L7").Value = 2900
$ws.Range("M7").Value = -2887.3333
$ws.Range("N7").Value = -3124

$ws.Range("H40").Value = 83846.47
$ws.Range("I40").Value = 172571.28
$ws.Range("K40").Value = 172571.28
$ws.Range("M40").Value = -172435.28

$ws.Range("H61").Value = 2905.2068
$ws.Range("J61").Value = 3466.1667
$ws.Range("L61").Value = 3466.1667
$ws.Range("N61").Value = -3870.1667

$ws.Range("H68").Value = 2041.7778
$ws.Range("I68").Value = 1655.1428
$ws.Range("J68").Value = 3395
$ws.Range("K68").Value = 1655.1428
$ws.Range("L68").Value = 3395
$ws.Range("M68").Value = -906.1428000000001
$ws.Range("N68").Value = -4893

$ws.Range("H71").Value = 2041.7778
$ws.Range("I71").Value = 1655.1428
$ws.Range("J71").Value = 3395
$ws.Range("K71").Value = 8275.714
$ws.Range("L71").Value = 16975
$ws.Range("M71").Value = -4531.714
$ws.Range("N71").Value = -24463

$ws.Range("H93").Value = 2494.2727
$ws.Range("I93").Value = 2034.4286
$ws.Range("K93").Value = 2034.4286
$ws.Range("M93").Value = -786.4286

$ws.Range("H113").Value = 2905.2068
$ws.Range("J113").Value = 3466.1667
$ws.Range("L113").Value = 3466.1667
$ws.Range("N113").Value = -7806.1667

$ws.Range("H122").Value = 4250.8096
$ws.Range("I122").Value = 4061.7856
$ws.Range("K122").Value = 12185.3568
$ws.Range("M122").Value = -9735.356800000001

$ws.Range("H126").Value = 2974.5
$ws.Range("I126").Value = 2999.3333
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 8997.999899999999
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -6527.999899999999
$ws.Range("N126").Value = -13640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9143.714
$ws.Range("J136").Value = 8201.200000000001
$ws.Range("L136").Value = 24603.6
$ws.Range("N136").Value = -29703.6
